$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$data = @(
    @("Danh mục", 18),
    @("Ngày công", 12),
    @("Phụ cấp", 420000),
    @("Lương cơ bản tại CẦN THƠ", 0),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Lương cơ bản tại LONG XUYÊN", 0),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Lương cơ bản tại SÓC TRĂNG", 0),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0)
)

$row = 1
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
